$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# Remove the Bharat Hardware & Co. entry (row 2) and its blank spacer row (row 3)
$ws1.Rows("2:3").Delete()

# Remove the trailing GST calculation block (old rows 36-40, now 34-38 after the shift above)
$ws1.Rows("34:38").Delete()

# Update Sale 22-23 F22 formula (additional 800000 deduction)
$ws2.Range("F22").Formula = "=E22-175496-500000-800000"

# Sheet view selections
$ws1.Application.ActiveWindow.ScrollRow = 19
$ws1.Range("D43").Select()
$ws2.Range("F23").Select()
